$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2621.2
$ws.Range("I62").Value = 1650
$ws.Range("J62").Value = 3268.6667
$ws.Range("K62").Value = 1650
$ws.Range("L62").Value = 3268.6667
$ws.Range("M62").Value = -1026
$ws.Range("N62").Value = -4516.6667
$ws.Range("H65").Value = 2621.2
$ws.Range("I65").Value = 1650
$ws.Range("J65").Value = 3268.6667
$ws.Range("K65").Value = 8250
$ws.Range("L65").Value = 16343.3335
$ws.Range("M65").Value = -5130
$ws.Range("N65").Value = -22583.3335
$ws.Range("H75").Value = 15313
$ws.Range("J75").Value = 15313
$ws.Range("L75").Value = 15313
$ws.Range("N75").Value = -17185
$ws.Range("H78").Value = 15313
$ws.Range("J78").Value = 15313
$ws.Range("L78").Value = 45939
$ws.Range("N78").Value = -55299
$ws.Range("H138").Value = 1185.0638
$ws.Range("I138").Value = 762.0755
$ws.Range("K138").Value = 2286.2265
$ws.Range("M138").Value = 2853.7735
$ws.Range("H141").Value = 1245.8334
$ws.Range("I141").Value = 895
$ws.Range("K141").Value = 2685
$ws.Range("M141").Value = 2495

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 10007.6
$ws.Range("J37").Value = 10007.6
$ws.Range("L37").Value = 10007.6
$ws.Range("N37").Value = -10553.6
$ws.Range("H61").Value = 1276.1666
$ws.Range("I61").Value = 1031.4
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 1031.4
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -819.4000000000001
$ws.Range("N61").Value = -2924
$ws.Range("H63").Value = 2237.5
$ws.Range("I63").Value = 2237.5
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2237.5
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1551.5
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2237.5
$ws.Range("I66").Value = 2237.5
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 11187.5
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -7755.5
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 1100.8823
$ws.Range("I74").Value = 572.8182
$ws.Range("J74").Value = 2069
$ws.Range("K74").Value = 572.8182
$ws.Range("L74").Value = 2069
$ws.Range("M74").Value = 301.1818
$ws.Range("N74").Value = -3817
$ws.Range("H77").Value = 1100.8823
$ws.Range("I77").Value = 572.8182
$ws.Range("J77").Value = 2069
$ws.Range("K77").Value = 2864.091
$ws.Range("L77").Value = 10345
$ws.Range("M77").Value = 1503.909
$ws.Range("N77").Value = -19081
$ws.Range("H96").Value = 18249.5
$ws.Range("J96").Value = 18249.5
$ws.Range("L96").Value = 18249.5
$ws.Range("N96").Value = -23741.5
$ws.Range("H136").Value = 1276.1666
$ws.Range("I136").Value = 1031.4
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 3094.2
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -544.2000000000003
$ws.Range("N136").Value = -12600

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1300
$ws.Range("I31").Value = 1300
$ws.Range("K31").Value = 1300
$ws.Range("M31").Value = -1005
$ws.Range("H34").Value = 1300
$ws.Range("I34").Value = 1300
$ws.Range("K34").Value = 1300
$ws.Range("M34").Value = -1098
$ws.Range("H92").Value = 57799.6
$ws.Range("J92").Value = 57799.6
$ws.Range("L92").Value = 57799.6
$ws.Range("N92").Value = -62791.6
$ws.Range("H96").Value = 14282.667
$ws.Range("J96").Value = 14282.667
$ws.Range("L96").Value = 14282.667
$ws.Range("N96").Value = -19774.667
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676
$ws.Range("H122").Value = 992.25
$ws.Range("I122").Value = 856.3333
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 2568.9999
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -118.9998999999998
$ws.Range("N122").Value = -9100
$ws.Range("H134").Value = 16668781
$ws.Range("I134").Value = 20835464
$ws.Range("J134").Value = 2053
$ws.Range("K134").Value = 62506392
$ws.Range("L134").Value = 6159
$ws.Range("M134").Value = -62503857
$ws.Range("N134").Value = -11229

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 8335095.5
$ws.Range("I34").Value = 574.5
$ws.Range("J34").Value = 10002000
$ws.Range("K34").Value = 1723.5
$ws.Range("L34").Value = 30006000
$ws.Range("M34").Value = -1639.5
$ws.Range("N34").Value = -30006168
$ws.Range("H39").Value = 1599.7693
$ws.Range("J39").Value = 1699.5
$ws.Range("L39").Value = 5098.5
$ws.Range("N39").Value = -5686.5
$ws.Range("H55").Value = 3285
$ws.Range("J55").Value = 3285
$ws.Range("L55").Value = 9855
$ws.Range("N55").Value = -10209
$ws.Range("H56").Value = 7810.769
$ws.Range("I56").Value = 7810.769
$ws.Range("K56").Value = 7810.769
$ws.Range("M56").Value = -7280.769
$ws.Range("H60").Value = 1543.5714
$ws.Range("I60").Value = 501.25
$ws.Range("J60").Value = 2933.3333
$ws.Range("K60").Value = 1503.75
$ws.Range("L60").Value = 8799.999899999999
$ws.Range("M60").Value = -1252.75
$ws.Range("N60").Value = -9301.999899999999
$ws.Range("H104").Value = 3752.611
$ws.Range("I104").Value = 4526
$ws.Range("J104").Value = 3655.9375
$ws.Range("K104").Value = 13578
$ws.Range("L104").Value = 10967.8125
$ws.Range("M104").Value = -10957
$ws.Range("N104").Value = -16209.8125
$ws.Range("H131").Value = 10205235
$ws.Range("J131").Value = 1216.6703
$ws.Range("L131").Value = 3650.0109
$ws.Range("N131").Value = -13730.0109
$ws.Range("H137").Value = 2671
$ws.Range("I137").Value = 966.36365
$ws.Range("J137").Value = 3657.8948
$ws.Range("K137").Value = 2899.09095
$ws.Range("L137").Value = 10973.6844
$ws.Range("M137").Value = 2200.90905
$ws.Range("N137").Value = -21173.6844

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 29997.8
$ws.Range("J92").Value = 29997.8
$ws.Range("L92").Value = 29997.8
$ws.Range("N92").Value = -33741.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 712.1177
$ws.Range("I16").Value = 820
$ws.Range("J16").Value = 514.3333
$ws.Range("K16").Value = 820
$ws.Range("L16").Value = 514.3333
$ws.Range("M16").Value = -650
$ws.Range("N16").Value = -854.3333
$ws.Range("H128").Value = 90000
$ws.Range("J128").Value = 90000
$ws.Range("L128").Value = 90000
$ws.Range("N128").Value = -99960
$ws.Range("H133").Value = 34019.8
$ws.Range("J133").Value = 34019.8
$ws.Range("L133").Value = 34019.8
$ws.Range("N133").Value = -39079.8
$ws.Range("H136").Value = 7664.9414
$ws.Range("I136").Value = 17881
$ws.Range("J136").Value = 2092.5454
$ws.Range("K136").Value = 53643
$ws.Range("L136").Value = 6277.6362
$ws.Range("M136").Value = -51093
$ws.Range("N136").Value = -11377.6362

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 1750
$ws.Range("J39").Value = 1750
$ws.Range("L39").Value = 1750
$ws.Range("N39").Value = -2576
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 14447479
$ws.Range("I122").Value = 16253164
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 48759492
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -48757042
$ws.Range("N122").Value = -10900

Write-Output "Applied all Kujata_Profits updates"